# Trade #44 closed at 2026-02-16 21:29:37 - leadlag DOWN +0.000%
# Append the new trade row (row 34) to the "leadlag" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("leadlag")

$row = 34

$ws.Cells.Item($row, 1).Value = 44             # Trade #

# B34 holds a date formatted as plain text ("2026-02-16"); Excel would
# otherwise auto-detect this as a date literal, so force a text format
# before writing it and drop back to the Normal style afterwards.
$dateCell = $ws.Cells.Item($row, 2)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2026-02-16"                 # Date
$dateCell.Style = "Normal"

$ws.Cells.Item($row, 3).Value = "21:29:37"                           # Time
$ws.Cells.Item($row, 4).Value = "leadlag"                            # Strategy
$ws.Cells.Item($row, 5).Value = "DOWN"                                # Side
$ws.Cells.Item($row, 6).Value = 68717.925                            # Entry Price
# G34 (Exit Price) stays blank - trade is still OPEN
$ws.Cells.Item($row, 8).Value = "OPEN"                                # Status
$ws.Cells.Item($row, 9).Value = 0                                    # P&L %
$ws.Cells.Item($row, 10).Value = 0                                   # P&L $
$ws.Cells.Item($row, 11).Value = 0.604                               # Confidence
$ws.Cells.Item($row, 12).Value = "Binance leading with -0.060% move" # Entry Reason
# M34 (Exit Reason) stays blank - trade is still OPEN
$ws.Cells.Item($row, 14).Value = 0                                   # Duration (min)
